$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B1: replace old Russian translation with the new one
$ws.Range("B1").Value = "9.5.2 Количество исследователей (в эквиваленте полной занятости) на миллион жителей"

# Add new column Q (2023) mirroring the formatting of column P
$ws.Range("P4").Copy($ws.Range("Q4"))
$ws.Range("Q4").Value = 2023

$ws.Range("P5").Copy($ws.Range("Q5"))
$ws.Range("Q5").Value = 631
